$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 21 (shifts old rows 21-32 down to 22-33)
$ws.Rows.Item(21).Insert()

# Fill in new row 21 for the BAT60J part covering D6
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = "BAT60J"
$ws.Cells.Item(21, 3).Value = "BAT60J"
$ws.Cells.Item(21, 4).Value = "SOD323"
$ws.Cells.Item(21, 5).Value = "D6"
$ws.Cells.Item(21, 6).Value = "SCHOTTKY DIODE"
$ws.Cells.Item(21, 7).Value = "STMicroelectronics"
$ws.Cells.Item(21, 8).Value = "BAT60JFILM"
$ws.Cells.Item(21, 9).Value = "497-3707-6-ND"
$ws.Cells.Item(21, 10).Value = "511-BAT60JFILM"
$ws.Cells.Item(21, 11).Value = "89K1218"

# Update row 20: quantity changes from 6 (covering D1-D6) to 5 (covering D1-D5)
$ws.Cells.Item(20, 1).Value = 5
$ws.Cells.Item(20, 5).Value = "D1, D2, D3, D4, D5"

# Update selection to reflect the active cell in the saved file
$ws.Range("B21").Select()
